$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "Supplier Product Number" in E1 (previously an empty,
# but styled, cell) -- this is inserted into the shared-string table right
# after "Short Description" and shifts all later string indices by one.
$ws.Range("E1").Value = "Supplier Product Number"

# Update the active selection left behind after editing, as Excel records it.
$ws.Range("H18").Select()
